# Apply updated market/profit data values to the Leve profit sheets.
# Source data: scheduled market data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 20908.555
$ws.Range("I43").Value = 30095.2
$ws.Range("K43").Value = 30095.2
$ws.Range("M43").Value = -30026.2

# Row 86
$ws.Range("H86").Value = 2242.1765
$ws.Range("I86").Value = 2185.3333
$ws.Range("K86").Value = 2185.3333
$ws.Range("M86").Value = -1062.3333

# Row 89
$ws.Range("H89").Value = 2242.1765
$ws.Range("I89").Value = 2185.3333
$ws.Range("K89").Value = 10926.6665
$ws.Range("M89").Value = -5310.666499999999

# Row 95
$ws.Range("H95").Value = 87999
$ws.Range("J95").Value = 87999
$ws.Range("L95").Value = 87999
$ws.Range("N95").Value = -93491

# Row 115
$ws.Range("H115").Value = 1136.9412
$ws.Range("I115").Value = 456
$ws.Range("K115").Value = 1368
$ws.Range("M115").Value = 199

# Row 127
$ws.Range("H127").Value = 3646.6191
$ws.Range("I127").Value = 3729.2
$ws.Range("J127").Value = 1995
$ws.Range("K127").Value = 11187.6
$ws.Range("L127").Value = 5985
$ws.Range("M127").Value = -6227.599999999999
$ws.Range("N127").Value = -15905

# Row 129
$ws.Range("H129").Value = 436962.62
$ws.Range("I129").Value = 465760.12
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 1397280.36
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -1392280.36
$ws.Range("N129").Value = -25000

# Row 132
$ws.Range("H132").Value = 11348.479
$ws.Range("I132").Value = 12418.35
$ws.Range("J132").Value = 4216
$ws.Range("K132").Value = 37255.05
$ws.Range("L132").Value = 12648
$ws.Range("M132").Value = -34725.05
$ws.Range("N132").Value = -17708

# Row 137
$ws.Range("H137").Value = 31258290
$ws.Range("I137").Value = 47621012
$ws.Range("K137").Value = 142863036
$ws.Range("M137").Value = -142860486

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 146143.27
$ws.Range("I32").Value = 150000.42
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 150000.42
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -149713.42
$ws.Range("N32").Value = -15574

# Row 61
$ws.Range("H61").Value = 4004340.5
$ws.Range("I61").Value = 4797.5454
$ws.Range("K61").Value = 4797.5454
$ws.Range("M61").Value = -4585.5454

# Row 62
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248

# Row 65
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240

# Row 74
$ws.Range("H74").Value = 2427118.5
$ws.Range("I74").Value = 4275823.5
$ws.Range("K74").Value = 4275823.5
$ws.Range("M74").Value = -4274949.5

# Row 77
$ws.Range("H77").Value = 2427118.5
$ws.Range("I77").Value = 4275823.5
$ws.Range("K77").Value = 21379117.5
$ws.Range("M77").Value = -21374749.5

# Row 110
$ws.Range("H110").Value = 1193.3846
$ws.Range("I110").Value = 1047.1818
$ws.Range("J110").Value = 1997.5
$ws.Range("K110").Value = 1047.1818
$ws.Range("L110").Value = 1997.5
$ws.Range("M110").Value = 997.8181999999999
$ws.Range("N110").Value = -6087.5

# Row 132
$ws.Range("H132").Value = 758903.4
$ws.Range("J132").Value = 1780.25
$ws.Range("L132").Value = 5340.75
$ws.Range("N132").Value = -10400.75

# Row 134
$ws.Range("H134").Value = 80141
$ws.Range("J134").Value = 80141
$ws.Range("L134").Value = 80141
$ws.Range("N134").Value = -90281

# Row 136
$ws.Range("H136").Value = 4004340.5
$ws.Range("I136").Value = 4797.5454
$ws.Range("K136").Value = 14392.6362
$ws.Range("M136").Value = -11842.6362

$ws = $wb.Worksheets.Item("BSM")
# Row 56
$ws.Range("H56").Value = 41766.668
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

# Row 80
$ws.Range("H80").Value = 1950
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 900
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = 98
$ws.Range("N80").Value = -4996

# Row 83
$ws.Range("H83").Value = 1950
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = 492
$ws.Range("N83").Value = -24984

# Row 105
$ws.Range("H105").Value = 5608.8184
$ws.Range("I105").Value = 3462.25
$ws.Range("J105").Value = 11333
$ws.Range("K105").Value = 3462.25
$ws.Range("L105").Value = 11333
$ws.Range("M105").Value = -1715.25
$ws.Range("N105").Value = -14827

# Row 132
$ws.Range("H132").Value = 194980
$ws.Range("J132").Value = 194980
$ws.Range("L132").Value = 194980
$ws.Range("N132").Value = -205100

# Row 134
$ws.Range("H134").Value = 16679858
$ws.Range("I134").Value = 4372.4287
$ws.Range("K134").Value = 13117.2861
$ws.Range("M134").Value = -10582.2861

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2528282.2
$ws.Range("J31").Value = 2871.9285
$ws.Range("L31").Value = 2871.9285
$ws.Range("N31").Value = -3461.9285

# Row 34
$ws.Range("H34").Value = 2528282.2
$ws.Range("J34").Value = 2871.9285
$ws.Range("L34").Value = 2871.9285
$ws.Range("N34").Value = -3275.9285

# Row 58
$ws.Range("H58").Value = 4912473.5
$ws.Range("I58").Value = 4958
$ws.Range("J58").Value = 10433428
$ws.Range("K58").Value = 4958
$ws.Range("L58").Value = 10433428
$ws.Range("M58").Value = -4755
$ws.Range("N58").Value = -10433834

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = ""
$ws.Range("N104").Value = 0

# Row 134
$ws.Range("H134").Value = 4147.5264
$ws.Range("I134").Value = 3058.3333
$ws.Range("J134").Value = 8232
$ws.Range("K134").Value = 9174.999899999999
$ws.Range("L134").Value = 24696
$ws.Range("M134").Value = -6639.999899999999
$ws.Range("N134").Value = -29766

# Row 136
$ws.Range("H136").Value = 4912473.5
$ws.Range("I136").Value = 4958
$ws.Range("J136").Value = 10433428
$ws.Range("K136").Value = 14874
$ws.Range("L136").Value = 31300284
$ws.Range("M136").Value = -12324
$ws.Range("N136").Value = -31305384

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5362206
$ws.Range("I5").Value = 5103149.5
$ws.Range("J5").Value = 5815554
$ws.Range("K5").Value = 15309448.5
$ws.Range("L5").Value = 17446662
$ws.Range("M5").Value = -15309336.5
$ws.Range("N5").Value = -17446886

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = ""
$ws.Range("M87").Value = ""
$ws.Range("N87").Value = 0

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = ""
$ws.Range("M90").Value = ""
$ws.Range("N90").Value = 0

# Row 109
$ws.Range("H109").Value = 1999.1428
$ws.Range("I109").Value = 1246.8
$ws.Range("K109").Value = 3740.4
$ws.Range("M109").Value = -2700.4

# Row 122
$ws.Range("H122").Value = 1076174.1
$ws.Range("I122").Value = 2016758.2
$ws.Range("J122").Value = 1220.7142
$ws.Range("K122").Value = 18150823.8
$ws.Range("L122").Value = 10986.4278
$ws.Range("M122").Value = -18148373.8
$ws.Range("N122").Value = -15886.4278

# Row 135
$ws.Range("H135").Value = 5362206
$ws.Range("I135").Value = 5103149.5
$ws.Range("J135").Value = 5815554
$ws.Range("K135").Value = 45928345.5
$ws.Range("L135").Value = 52339986
$ws.Range("M135").Value = -45925810.5
$ws.Range("N135").Value = -52345056

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 16252
$ws.Range("J10").Value = 16252
$ws.Range("L10").Value = 16252
$ws.Range("N10").Value = -16590

# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = ""
$ws.Range("N59").Value = 0

# Row 70
$ws.Range("H70").Value = 45782.145
$ws.Range("I70").Value = 30462.334
$ws.Range("K70").Value = 30462.334
$ws.Range("M70").Value = -30192.334

# Row 73
$ws.Range("H73").Value = 45782.145
$ws.Range("I73").Value = 30462.334
$ws.Range("K73").Value = 30462.334
$ws.Range("M73").Value = -29526.334

# Row 102
$ws.Range("H102").Value = 3125.95
$ws.Range("I102").Value = 2923.2222
$ws.Range("J102").Value = 4950.5
$ws.Range("K102").Value = 2923.2222
$ws.Range("L102").Value = 4950.5
$ws.Range("M102").Value = -1301.2222
$ws.Range("N102").Value = -8194.5

# Row 122
$ws.Range("H122").Value = 3744.7778
$ws.Range("I122").Value = 901.3333
$ws.Range("J122").Value = 5166.5
$ws.Range("K122").Value = 2703.9999
$ws.Range("L122").Value = 15499.5
$ws.Range("M122").Value = -253.9998999999998
$ws.Range("N122").Value = -20399.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3039.1333
$ws.Range("I22").Value = 2473.75
$ws.Range("J22").Value = 3244.7273
$ws.Range("K22").Value = 2473.75
$ws.Range("L22").Value = 3244.7273
$ws.Range("M22").Value = -2178.75
$ws.Range("N22").Value = -3834.7273

# Row 27
$ws.Range("H27").Value = 3039.1333
$ws.Range("I27").Value = 2473.75
$ws.Range("J27").Value = 3244.7273
$ws.Range("K27").Value = 2473.75
$ws.Range("L27").Value = 3244.7273
$ws.Range("M27").Value = -2366.75
$ws.Range("N27").Value = -3458.7273

# Row 55
$ws.Range("H55").Value = 1525.7391
$ws.Range("J55").Value = 1385.4286
$ws.Range("L55").Value = 1385.4286
$ws.Range("N55").Value = -1731.4286

# Row 122
$ws.Range("H122").Value = 3565.6
$ws.Range("I122").Value = 3186.5
$ws.Range("J122").Value = 3998.8572
$ws.Range("K122").Value = 9559.5
$ws.Range("L122").Value = 11996.5716
$ws.Range("M122").Value = -7109.5
$ws.Range("N122").Value = -16896.5716

# Row 136
$ws.Range("H136").Value = 10420352
$ws.Range("I136").Value = 11367786
$ws.Range("J136").Value = 9618678
$ws.Range("K136").Value = 34103358
$ws.Range("L136").Value = 28856034
$ws.Range("M136").Value = -34100808
$ws.Range("N136").Value = -28861134
